$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 2's Email/Password/Expected-Result values: reuse the same
# Email/Password pairs used elsewhere (TC_02's email, TC_03's password)
# and drop the two now-unused strings (the old unique email & password).
$ws.Range("C2").Value = "CorrectEmail@gmail.com"
$ws.Range("D2").Value = "thatshouldmakeit"

# Make C2 a mailto hyperlink, same as C3/C4.
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:CorrectEmail@gmail.com") | Out-Null

# Hyperlinks.Add() minted a fresh (duplicate) cell style for C2; restore
# the original shared "Hyperlink" style (same one C3/C4 already use) by
# copying formats over from C3.
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)

# Restore the selection to D2 (was G7).
$ws.Range("D2").Select() | Out-Null
